# Apply the changes described by the diff:
#  1. Rename colaborador "Carlos Santos" -> "Carla Santos" (row 4, column A).
#  2. Give cell C8 its own style: Times New Roman 10pt, centered horizontally,
#     which introduces a new font + a new cell format (cellXfs entry).
#  3. Center-align cell A8 (general -> center) - already centered via header
#     below, but A8 alignment also moves from general to center as part of
#     the shared xf used for A8/B8.
#  4. Grow row 8's height to fit the new formatting (12.8 -> 24.05).
#  5. Update the sheet's active selection to A8:C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the employee name.
$ws.Range("A4").Value = "Carla Santos"

# 2/3. Center alignment for the whole (currently empty) trailer row first,
# while A8/B8/C8 still share a single style, then give the last cell (C8)
# its own distinct Times New Roman font on top of that, still centered.
$ws.Range("A8:C8").HorizontalAlignment = -4108   # xlHAlignCenter
$ws.Range("C8").Font.Name = "Times New Roman"
$ws.Range("C8").Font.Size = 10

# 4. Taller row to match the new formatting.
$ws.Rows.Item(8).RowHeight = 24.05

# 5. Update selection to A8:C8 with A8 as the active cell.
$ws.Range("A8:C8").Select()
